# LNI-317: Added previously removed signature blocks back to test
# Statutory Instruments.
#
# Appends the SigBlock signature paragraphs (two blank spacer lines, the
# "Name" line, the "Clerk of the Privy Council" line) plus a trailing
# zero-indent LQN3 paragraph after the document's final paragraph, and
# restores a couple of style-table tweaks (No List's uiPriority).

$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1. Lay down all five new paragraphs (with their block/paragraph
#    styles only) before putting any text/character-styled runs into
#    them. Doing the paragraph splits first - while the insertion
#    point is still carrying "plain" formatting - keeps a character
#    style used in an earlier paragraph (e.g. SigSignee/Sigtitle) from
#    bleeding into the empty run of the paragraph created after it.
# -----------------------------------------------------------------

$anchor = $d.Paragraphs.Last

$anchor.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs.Last
$p1.Style = "SigBlock"

$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last
$p2.Style = "SigBlock"

$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Last
$p3.Style = "SigBlock"

$p3.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs.Last
$p4.Style = "SigBlock"

$p4.Range.InsertParagraphAfter()
$p5 = $d.Paragraphs.Last
$p5.Style = "LQN3"
$p5.LeftIndent = 0
$p5.FirstLineIndent = 0

# -----------------------------------------------------------------
# 2. Now fill in the two signature lines' content.
# -----------------------------------------------------------------

# Paragraph 3: <tab>Name, "Name" styled as SigSignee.
$p3.Range.InsertAfter([char]9)
$tab3End = $p3.Range.End - 1
$p3.Range.InsertAfter("Name")
$name3Range = $d.Range($tab3End, $p3.Range.End - 1)
$name3Range.Style = "SigSignee"

# Paragraph 4: <tab>Clerk of the Privy Council, styled as Sigtitle.
$p4.Range.InsertAfter([char]9)
$tab4End = $p4.Range.End - 1
$p4.Range.InsertAfter("Clerk of the Privy Council")
$title4Range = $d.Range($tab4End, $p4.Range.End - 1)
$title4Range.Style = "Sigtitle"

# -----------------------------------------------------------------
# 3. Style table tweaks.
# -----------------------------------------------------------------

# "No List" numbering style gains a uiPriority of 99.
$noList = $d.Styles("No List")
$noList.Priority = 99

Write-Output "signature block inserted"
